$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 4.394810730541869
$ws.Range("D2").Value = 8.705329417287764
$ws.Range("E2").Value = 13.48940174833411
$ws.Range("F2").Value = 35.20292828826183
$ws.Range("G2").Value = 38.57583020164487
$ws.Range("H2").Value = 16.69741828078335
$ws.Range("I2").Value = 24.51744850372801
$ws.Range("J2").Value = 10.01631730254438
$ws.Range("K2").Value = 18.203195154462
$ws.Range("C3").Value = 4.333750065981086
$ws.Range("D3").Value = 8.637066540057091
$ws.Range("E3").Value = 13.42484874682014
$ws.Range("F3").Value = 35.3145944369608
$ws.Range("G3").Value = 38.73748296874551
$ws.Range("H3").Value = 16.78530786875513
$ws.Range("I3").Value = 24.62476834721426
$ws.Range("J3").Value = 10.02738553507294
$ws.Range("K3").Value = 17.59373171871528
$ws.Range("C4").Value = 4.29558329770927
$ws.Range("D4").Value = 8.595718302031672
$ws.Range("E4").Value = 13.38766599454569
$ws.Range("F4").Value = 35.3958590508314
$ws.Range("G4").Value = 38.85553599774195
$ws.Range("H4").Value = 16.84368780417193
$ws.Range("I4").Value = 24.69922376114399
$ws.Range("J4").Value = 10.03596043459937
$ws.Range("K4").Value = 17.21013151882299
$ws.Range("C5").Value = 4.279871241662373
$ws.Range("D5").Value = 8.579023086657886
$ws.Range("E5").Value = 13.37314129313603
$ws.Range("F5").Value = 35.4321478298917
$ws.Range("G5").Value = 38.90832365574826
$ws.Range("H5").Value = 16.86858450752073
$ws.Range("I5").Value = 24.7317033375167
$ws.Range("J5").Value = 10.03990174011817
$ws.Range("K5").Value = 17.05167497359509
$ws.Range("C6").Value = 4.277253003285816
$ws.Range("D6").Value = 8.576260547435833
$ws.Range("E6").Value = 13.37076767876887
$ws.Range("F6").Value = 35.43836454602751
$ws.Range("G6").Value = 38.91737026448045
$ws.Range("H6").Value = 16.87278529774835
$ws.Range("I6").Value = 24.73722527656233
$ws.Range("J6").Value = 10.04058317354736
$ws.Range("K6").Value = 17.02524129645297
$ws.Range("C7").Value = 4.295372027623211
$ws.Range("D7").Value = 8.595492502718695
$ws.Range("E7").Value = 13.38746755475226
$ws.Range("F7").Value = 35.39633563625741
$ws.Range("G7").Value = 38.8562290245865
$ws.Range("H7").Value = 16.84401909518886
$ws.Range("I7").Value = 24.69965315299981
$ws.Range("J7").Value = 10.03601177925608
$ws.Range("K7").Value = 17.20800285761946
$ws.Range("C8").Value = 4.373903687876924
$ws.Range("D8").Value = 8.681682495844006
$ws.Range("E8").Value = 13.46664142450758
$ws.Range("F8").Value = 35.23878236105647
$ws.Range("G8").Value = 38.62763922435028
$ws.Range("H8").Value = 16.72680386130776
$ws.Range("I8").Value = 24.55266718206595
$ws.Range("J8").Value = 10.01976421299852
$ws.Range("K8").Value = 17.99512876317281
$ws.Range("C9").Value = 4.522077882782558
$ws.Range("D9").Value = 8.854635883426337
$ws.Range("E9").Value = 13.64085767410134
$ws.Range("F9").Value = 35.03148386490798
$ws.Range("G9").Value = 38.33049933071435
$ws.Range("H9").Value = 16.53215368048813
$ws.Range("I9").Value = 24.33297190285307
$ws.Range("J9").Value = 10.0020313148239
$ws.Range("K9").Value = 19.45536189102431
$ws.Range("C10").Value = 4.626787328268454
$ws.Range("D10").Value = 8.983355186950124
$ws.Range("E10").Value = 13.77968266853624
$ws.Range("F10").Value = 34.9422837388978
$ws.Range("G10").Value = 38.2068347053793
$ws.Range("H10").Value = 16.41084913379249
$ws.Range("I10").Value = 24.21416717901758
$ws.Range("J10").Value = 9.997629409782729
$ws.Range("K10").Value = 20.46694930618154
$ws.Range("C11").Value = 4.673400040147295
$ws.Range("D11").Value = 9.042110080685216
$ws.Range("E11").Value = 13.84501984184824
$ws.Range("F11").Value = 34.91561760910075
$ws.Range("G11").Value = 38.17160357233822
$ws.Range("H11").Value = 16.36043022498756
$ws.Range("I11").Value = 24.16953659878873
$ws.Range("J11").Value = 9.997500559179228
$ws.Range("K11").Value = 20.91201179897362
$ws.Range("C12").Value = 4.690894801223134
$ws.Range("D12").Value = 9.064374130415146
$ws.Range("E12").Value = 13.87006032540015
$ws.Range("F12").Value = 34.90753514163497
$ws.Range("G12").Value = 38.1613197812197
$ws.Range("H12").Value = 16.34202731219048
$ws.Range("I12").Value = 24.15400265838429
$ws.Range("J12").Value = 9.997720985059031
$ws.Range("K12").Value = 21.07824207071956
$ws.Range("C13").Value = 4.68713410086879
$ws.Range("D13").Value = 9.059578729131863
$ws.Range("E13").Value = 13.86465437570377
$ws.Range("F13").Value = 34.90918598909766
$ws.Range("G13").Value = 38.16339809250346
$ws.Range("H13").Value = 16.34595996180125
$ws.Range("I13").Value = 24.15728718724055
$ws.Range("J13").Value = 9.997661543422488
$ws.Range("K13").Value = 21.04254586631055
$ws.Range("C14").Value = 4.674842534760437
$ws.Range("D14").Value = 9.043941539244742
$ws.Range("E14").Value = 13.84707404129681
$ws.Range("F14").Value = 34.91491218797417
$ws.Range("G14").Value = 38.17069605959562
$ws.Range("H14").Value = 16.35890235792803
$ws.Range("I14").Value = 24.16823114981521
$ws.Range("J14").Value = 9.997513299822463
$ws.Range("K14").Value = 20.9257345632583
$ws.Range("C15").Value = 4.667292937392601
$ws.Range("D15").Value = 9.03436483118421
$ws.Range("E15").Value = 13.83634401208462
$ws.Range("F15").Value = 34.91868253226363
$ws.Range("G15").Value = 38.17556540044444
$ws.Range("H15").Value = 16.36691990010161
$ws.Range("I15").Value = 24.17511299947541
$ws.Range("J15").Value = 9.997457548359243
$ws.Range("K15").Value = 20.85388041308184
$ws.Range("C16").Value = 4.623719705240656
$ws.Range("D16").Value = 8.979518249188734
$ws.Range("E16").Value = 13.77545538621671
$ws.Range("F16").Value = 34.9443076849884
$ws.Range("G16").Value = 38.20956321388029
$ws.Range("H16").Value = 16.414240375434
$ws.Range("I16").Value = 24.21727445127449
$ws.Range("J16").Value = 9.997675517711768
$ws.Range("K16").Value = 20.43754765892592
$ws.Range("C17").Value = 4.596720569134023
$ws.Range("D17").Value = 8.945912933774006
$ws.Range("E17").Value = 13.73865096282277
$ws.Range("F17").Value = 34.96360212265508
$ws.Range("G17").Value = 38.23583011202752
$ws.Range("H17").Value = 16.44449323942566
$ws.Range("I17").Value = 24.24555975501623
$ws.Range("J17").Value = 9.998289005251095
$ws.Range("K17").Value = 20.17817105391748
$ws.Range("C18").Value = 4.581095916591138
$ws.Range("D18").Value = 8.9266036843187
$ws.Range("E18").Value = 13.71768865497234
$ws.Range("F18").Value = 34.97600839403632
$ws.Range("G18").Value = 38.25291524351254
$ws.Range("H18").Value = 16.46234193213406
$ws.Range("I18").Value = 24.26271409547006
$ws.Range("J18").Value = 9.998818239281862
$ws.Range("K18").Value = 20.02756874373302
$ws.Range("C19").Value = 4.575789586490748
$ws.Range("D19").Value = 8.920069701282383
$ws.Range("E19").Value = 13.7106271287166
$ws.Range("F19").Value = 34.98043321767733
$ws.Range("G19").Value = 38.25903844484395
$ws.Range("H19").Value = 16.46846200243274
$ws.Range("I19").Value = 24.26867392139426
$ws.Range("J19").Value = 9.999027724964076
$ws.Range("K19").Value = 19.97633851532866
$ws.Range("C20").Value = 4.599604625521154
$ws.Range("D20").Value = 8.949488345153808
$ws.Range("E20").Value = 13.74254758236326
$ws.Range("F20").Value = 34.96141264610448
$ws.Range("G20").Value = 38.23282907430192
$ws.Range("H20").Value = 16.44122636564795
$ws.Range("I20").Value = 24.24245700711022
$ws.Range("J20").Value = 9.998205445682695
$ws.Range("K20").Value = 20.20592966680412
$ws.Range("C21").Value = 4.678457187638353
$ws.Range("D21").Value = 9.048534268421163
$ws.Range("E21").Value = 13.85222983053777
$ws.Range("F21").Value = 34.91317545582614
$ws.Range("G21").Value = 38.16846923657521
$ws.Range("H21").Value = 16.35508210833479
$ws.Range("I21").Value = 24.16497945118262
$ws.Range("J21").Value = 9.997549538294113
$ws.Range("K21").Value = 20.9601084108977
$ws.Range("C22").Value = 4.729075702950357
$ws.Range("D22").Value = 9.113346654053817
$ws.Range("E22").Value = 13.9256471338539
$ws.Range("F22").Value = 34.8934027893843
$ws.Range("G22").Value = 38.14424088598898
$ws.Range("H22").Value = 16.30280377798276
$ws.Range("I22").Value = 24.12231505813051
$ws.Range("J22").Value = 9.998689892621359
$ws.Range("K22").Value = 21.43952083724862
$ws.Range("C23").Value = 4.702146592261505
$ws.Range("D23").Value = 9.078752300646634
$ws.Range("E23").Value = 13.8863095528951
$ws.Range("F23").Value = 34.90287591833706
$ws.Range("G23").Value = 38.15552981887619
$ws.Range("H23").Value = 16.33033613766904
$ws.Range("I23").Value = 24.14435233320609
$ws.Range("J23").Value = 9.997937798266934
$ws.Range("K23").Value = 21.18492301745649
$ws.Range("C24").Value = 4.598301062111743
$ws.Range("D24").Value = 8.947871866695269
$ws.Range("E24").Value = 13.74078530575128
$ws.Range("F24").Value = 34.96239841772989
$ws.Range("G24").Value = 38.23417966631118
$ws.Range("H24").Value = 16.44270189947219
$ws.Range("I24").Value = 24.24385697960464
$ws.Range("J24").Value = 9.998242673079659
$ws.Range("K24").Value = 20.19338460772255
$ws.Range("C25").Value = 4.482680141465091
$ws.Range("D25").Value = 8.807502895999537
$ws.Range("E25").Value = 13.59177047381946
$ws.Range("F25").Value = 35.07655818377351
$ws.Range("G25").Value = 38.39443390602783
$ws.Range("H25").Value = 16.58101769954775
$ws.Range("I25").Value = 24.38498077103814
$ws.Range("J25").Value = 10.00531412955262
$ws.Range("K25").Value = 19.07040310500382
